# Peru Liga 1 - base update (19-04-2024 00:38)
# Re-applies the reordering of a handful of match rows (the underlying
# source re-sorted some rows) plus refreshed odds for the last upcoming
# fixture and a newly added fixture.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Simple 2-row swaps (row content in columns B:AC swaps between the
#    two rows; column A - the sequential id - stays put).
# ---------------------------------------------------------------------
$swapPairs = @(
    @(61, 62),
    @(156, 157),
    @(228, 229),
    @(252, 253)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $rng1 = $ws.Range("B$r1`:AC$r1")
    $rng2 = $ws.Range("B$r2`:AC$r2")
    $v1 = $rng1.Value2
    $v2 = $rng2.Value2
    $rng1.Value2 = $v2
    $rng2.Value2 = $v1
}

# ---------------------------------------------------------------------
# 2) A 9-row reshuffle (rows 180-188): each row's content in columns
#    B:AC is replaced by another row's original content (column A keeps
#    its sequential value).
# ---------------------------------------------------------------------
$blockRows = 180, 181, 182, 183, 184, 185, 186, 187, 188
$snapshot = @{}
foreach ($r in $blockRows) {
    $snapshot[$r] = $ws.Range("B$r`:AC$r").Value2
}

# after-row -> source before-row (content provider)
$permMap = @{
    180 = 182
    181 = 180
    182 = 181
    183 = 186
    184 = 185
    185 = 183
    186 = 184
    187 = 188
    188 = 187
}

foreach ($r in $blockRows) {
    $src = $permMap[$r]
    $ws.Range("B$r`:AC$r").Value2 = $snapshot[$src]
}

# ---------------------------------------------------------------------
# 3) Row 290: odds/metadata refreshed for the still-unplayed fixture
#    (result columns H/I/J/AB/AC remain blank).
# ---------------------------------------------------------------------
$ws.Range("B290").Value2 = 8042077
$ws.Range("E290").Value2 = 45401.70833333334
$ws.Range("F290").Value2 = "Union Comercio"
$ws.Range("G290").Value2 = "Universitario de Deportes"
$ws.Range("K290").Value2 = 5
$ws.Range("L290").Value2 = 4.333
$ws.Range("M290").Value2 = 1.533
$ws.Range("N290").Value2 = 5.75
$ws.Range("O290").Value2 = 4.5
$ws.Range("P290").Value2 = 1.444
$ws.Range("Q290").Value2 = 1
$ws.Range("R290").Value2 = 2.1
$ws.Range("S290").Value2 = 1.775
$ws.Range("U290").Value2 = 1.975
$ws.Range("V290").Value2 = 1.875

# ---------------------------------------------------------------------
# 4) New row 291: a newly added upcoming fixture.
#    Copy formatting (bold/border style for A, date number-format for E)
#    from row 290 so no new cell styles are introduced.
# ---------------------------------------------------------------------
$ws.Range("A290").Copy() | Out-Null
$ws.Range("A291").PasteSpecial(-4122) | Out-Null
$ws.Range("E290").Copy() | Out-Null
$ws.Range("E291").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A291").Value2 = 289
$ws.Range("B291").Value2 = 8042078
$ws.Range("C291").Value2 = "Peru Liga 1"
$ws.Range("D291").Value2 = "Peru Liga 1"
$ws.Range("E291").Value2 = 45401.91666666666
$ws.Range("F291").Value2 = "Cienciano"
$ws.Range("G291").Value2 = "Deportivo Garcilaso"
$ws.Range("K291").Value2 = 1.727
$ws.Range("L291").Value2 = 3.75
$ws.Range("M291").Value2 = 4.333
$ws.Range("N291").Value2 = 2.15
$ws.Range("O291").Value2 = 3.4
$ws.Range("P291").Value2 = 3.4
$ws.Range("Q291").Value2 = -0.25
$ws.Range("R291").Value2 = 1.825
$ws.Range("S291").Value2 = 2.025
$ws.Range("T291").Value2 = 2.5
$ws.Range("U291").Value2 = 2.025
$ws.Range("V291").Value2 = 1.825
$ws.Range("W291").Value2 = 0
$ws.Range("X291").Value2 = 0
$ws.Range("Y291").Value2 = 0
$ws.Range("Z291").Value2 = 0
$ws.Range("AA291").Value2 = 0
